# Natmi following Dr Hou advice
# Update the NATMI ligand-receptor edge statistics (rows 2-10) to reflect
# the revised ligand/receptor-expressing cell counts (1 -> 3) and the
# recomputed expression / specificity values that follow from that change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 26.29132066666667
$ws.Cells.Item(2, 8).Value = 78.873962
$ws.Cells.Item(2, 9).Value = 0.1411782207947891
$ws.Cells.Item(2, 10).Value = 0.1411782207947891
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 36.51516
$ws.Cells.Item(2, 14).Value = 109.54548
$ws.Cells.Item(2, 15).Value = 0.3318909895952502
$ws.Cells.Item(2, 16).Value = 0.3318909895952502
$ws.Cells.Item(2, 17).Value = 960.0317807546401
$ws.Cells.Item(2, 18).Value = 8640.28602679176
$ws.Cells.Item(2, 19).Value = 0.04685577940887926
$ws.Cells.Item(2, 20).Value = 0.04685577940887927
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 26.29132066666667
$ws.Cells.Item(3, 8).Value = 78.873962
$ws.Cells.Item(3, 9).Value = 0.1411782207947891
$ws.Cells.Item(3, 10).Value = 0.1411782207947891
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 42.26455300000001
$ws.Cells.Item(3, 14).Value = 126.793659
$ws.Cells.Item(3, 15).Value = 0.3841479626536184
$ws.Cells.Item(3, 16).Value = 0.3841479626536184
$ws.Cells.Item(3, 17).Value = 1111.190915756329
$ws.Cells.Item(3, 18).Value = 10000.71824180696
$ws.Cells.Item(3, 19).Value = 0.05423332588938092
$ws.Cells.Item(3, 20).Value = 0.05423332588938093
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 26.29132066666667
$ws.Cells.Item(4, 8).Value = 78.873962
$ws.Cells.Item(4, 9).Value = 0.1411782207947891
$ws.Cells.Item(4, 10).Value = 0.1411782207947891
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 31.24183366666666
$ws.Cells.Item(4, 14).Value = 93.725501
$ws.Cells.Item(4, 15).Value = 0.2839610477511313
$ws.Cells.Item(4, 16).Value = 0.2839610477511314
$ws.Cells.Item(4, 17).Value = 821.3890671449958
$ws.Cells.Item(4, 18).Value = 7392.501604304962
$ws.Cells.Item(4, 19).Value = 0.04008911549652886
$ws.Cells.Item(4, 20).Value = 0.04008911549652887
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 143.4723713333333
$ws.Cells.Item(5, 8).Value = 430.417114
$ws.Cells.Item(5, 9).Value = 0.7704129577533824
$ws.Cells.Item(5, 10).Value = 0.7704129577533824
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 36.51516
$ws.Cells.Item(5, 14).Value = 109.54548
$ws.Cells.Item(5, 15).Value = 0.3318909895952502
$ws.Cells.Item(5, 16).Value = 0.3318909895952502
$ws.Cells.Item(5, 17).Value = 5238.91659481608
$ws.Cells.Item(5, 18).Value = 47150.24935334471
$ws.Cells.Item(5, 19).Value = 0.2556931189457737
$ws.Cells.Item(5, 20).Value = 0.2556931189457738
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 143.4723713333333
$ws.Cells.Item(6, 8).Value = 430.417114
$ws.Cells.Item(6, 9).Value = 0.7704129577533824
$ws.Cells.Item(6, 10).Value = 0.7704129577533824
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 42.26455300000001
$ws.Cells.Item(6, 14).Value = 126.793659
$ws.Cells.Item(6, 15).Value = 0.3841479626536184
$ws.Cells.Item(6, 16).Value = 0.3841479626536184
$ws.Cells.Item(6, 17).Value = 6063.795642253347
$ws.Cells.Item(6, 18).Value = 54574.16078028013
$ws.Cells.Item(6, 19).Value = 0.29595256812291
$ws.Cells.Item(6, 20).Value = 0.29595256812291
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 143.4723713333333
$ws.Cells.Item(7, 8).Value = 430.417114
$ws.Cells.Item(7, 9).Value = 0.7704129577533824
$ws.Cells.Item(7, 10).Value = 0.7704129577533824
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 31.24183366666666
$ws.Cells.Item(7, 14).Value = 93.725501
$ws.Cells.Item(7, 15).Value = 0.2839610477511313
$ws.Cells.Item(7, 16).Value = 0.2839610477511314
$ws.Cells.Item(7, 17).Value = 4482.339960958234
$ws.Cells.Item(7, 18).Value = 40341.05964862411
$ws.Cells.Item(7, 19).Value = 0.2187672706846985
$ws.Cells.Item(7, 20).Value = 0.2187672706846986
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 16.46418733333334
$ws.Cells.Item(8, 8).Value = 49.39256200000001
$ws.Cells.Item(8, 9).Value = 0.08840882145182853
$ws.Cells.Item(8, 10).Value = 0.08840882145182853
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 36.51516
$ws.Cells.Item(8, 14).Value = 109.54548
$ws.Cells.Item(8, 15).Value = 0.3318909895952502
$ws.Cells.Item(8, 16).Value = 0.3318909895952502
$ws.Cells.Item(8, 17).Value = 601.1924347466401
$ws.Cells.Item(8, 18).Value = 5410.731912719761
$ws.Cells.Item(8, 19).Value = 0.02934209124059715
$ws.Cells.Item(8, 20).Value = 0.02934209124059716
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 16.46418733333334
$ws.Cells.Item(9, 8).Value = 49.39256200000001
$ws.Cells.Item(9, 9).Value = 0.08840882145182853
$ws.Cells.Item(9, 10).Value = 0.08840882145182853
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 42.26455300000001
$ws.Cells.Item(9, 14).Value = 126.793659
$ws.Cells.Item(9, 15).Value = 0.3841479626536184
$ws.Cells.Item(9, 16).Value = 0.3841479626536184
$ws.Cells.Item(9, 17).Value = 695.8515181515955
$ws.Cells.Item(9, 18).Value = 6262.663663364359
$ws.Cells.Item(9, 19).Value = 0.03396206864132744
$ws.Cells.Item(9, 20).Value = 0.03396206864132745
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 16.46418733333334
$ws.Cells.Item(10, 8).Value = 49.39256200000001
$ws.Cells.Item(10, 9).Value = 0.08840882145182853
$ws.Cells.Item(10, 10).Value = 0.08840882145182853
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 31.24183366666666
$ws.Cells.Item(10, 14).Value = 93.725501
$ws.Cells.Item(10, 15).Value = 0.2839610477511313
$ws.Cells.Item(10, 16).Value = 0.2839610477511314
$ws.Cells.Item(10, 17).Value = 514.3714021248402
$ws.Cells.Item(10, 18).Value = 4629.342619123562
$ws.Cells.Item(10, 19).Value = 0.02510466156990393
$ws.Cells.Item(10, 20).Value = 0.02510466156990393
